$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the new columns, copying the existing header style (from AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record values for every data row (2 through 44)
$ws.Range("AD2:AD44").Value = 98
$ws.Range("AE2:AE44").Value = 64
$ws.Range("AF2:AF44").Value = 0
